# msz - mandatory field checks part 3 incl. groupbox
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header xpath locators. D1 switches from the "section[2]" xpath
#     to the new "enterinsurantdata" groupbox locator. (introduces new
#     shared string #46) ---
$ws.Cells.Item(1, 4).Value = '//*[@id="enterinsurantdata"]'

# --- New rows 20-22: Insurant Page mandatory-field checks (mirrors the
#     VehicleData block in rows 17-19, but targeting column D). Text is
#     populated in row 21, 22, 20 order to match the shared-string layout
#     of the authored workbook (new strings #47-#49). ---
$ws.Cells.Item(21, 1).Value = "102_AutomobileInsurance_003_InsurantData_001_MandatoryFields_FillFirstName"
$ws.Cells.Item(22, 1).Value = "102_AutomobileInsurance_003_InsurantData_001_MandatoryFields_CheckFilledFirstName"
$ws.Cells.Item(20, 1).Value = "Insurant Page check for open mandatory fields"

$ws.Cells.Item(20, 2).Value = "<CHK>"
$ws.Cells.Item(20, 4).Value = "Insurant Page check for open mandatory fields"
$ws.Cells.Item(20, 8).Value = "<NOP>"

$ws.Cells.Item(21, 2).Value = "<SET>"
$ws.Cells.Item(21, 4).Value = "102_AutomobileInsurance_003_InsurantData_001_MandatoryFields_FillFirstName"
$ws.Cells.Item(21, 8).Value = "<NOP>"

$ws.Cells.Item(22, 2).Value = "<CHK>"
$ws.Cells.Item(22, 4).Value = "102_AutomobileInsurance_003_InsurantData_001_MandatoryFields_CheckFilledFirstName"
$ws.Cells.Item(22, 8).Value = "<NOP>"

# --- Row 2: D2 gets the old section[2] locator, highlighted with the same
#     amber fill used by the other per-row override cells (A7:A10). ---
$ws.Cells.Item(2, 4).Value = '//*[@id="insurance-form"]/div/section[2]'
$ws.Cells.Item(2, 4).Interior.Color = 49407

# --- Column widths: column A widened (no longer auto "best fit"), and
#     column D split off from the D:E best-fit pair with its own width. ---
$ws.Columns.Item(1).ColumnWidth = 76.72135416666667
$ws.Columns.Item(4).ColumnWidth = 75.38541666666667

# --- Selection / active cell moved. ---
$ws.Range("D16").Select()

# --- Picture: shifted down by one row, narrowed; height unchanged. ---
$shp = $ws.Shapes.Item(1)
$shp.Top = 352.2
$shp.Width = 639.2079484498031
